$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4849909108651786
$ws.Range("C2").Value = 0.1390367536166366
$ws.Range("D2").Value = 0.01388781897824742
$ws.Range("F2").Value = 0.4548297057558699
$ws.Range("G2").Value = 0.002376288671846065
$ws.Range("I2").Value = 0.3201922625708562
$ws.Range("M2").Value = 0.8950069012951758
$ws.Range("N2").Value = 0.9088538593433313
$ws.Range("O2").Value = 1.432740183228361
$ws.Range("B3").Value = 0.4235675016159064
$ws.Range("C3").Value = 0.1243917632931186
$ws.Range("D3").Value = 0.0122904808623332
$ws.Range("F3").Value = 0.446652616678854
$ws.Range("G3").Value = 0.002378783334122204
$ws.Range("I3").Value = 0.3223965144529828
$ws.Range("M3").Value = 0.7963853923832573
$ws.Range("N3").Value = 0.9210814647999364
$ws.Range("O3").Value = 1.418427280022513
$ws.Range("B4").Value = 0.3857360183631329
$ws.Range("C4").Value = 0.1153398086294999
$ws.Range("D4").Value = 0.01130419820463402
$ws.Range("F4").Value = 0.441978754326982
$ws.Range("G4").Value = 0.002380397232916552
$ws.Range("I4").Value = 0.323991918952057
$ws.Range("M4").Value = 0.73634678032316
$ws.Range("N4").Value = 0.9290192930965944
$ws.Range("O4").Value = 1.410746336141926
$ws.Range("B5").Value = 0.3702911067308605
$ws.Range("C5").Value = 0.1116362770502519
$ws.Range("D5").Value = 0.01090092363537565
$ws.Range("F5").Value = 0.44016120513934
$ws.Range("G5").Value = 0.002381075635677651
$ws.Range("I5").Value = 0.324702764161767
$ws.Range("M5").Value = 0.7120042915382783
$ws.Range("N5").Value = 0.9323621219573575
$ws.Range("O5").Value = 1.407894133364493
$ws.Range("B6").Value = 0.3677248188329258
$ws.Range("C6").Value = 0.1110204232792711
$ws.Range("D6").Value = 0.01083387914005129
$ws.Range("F6").Value = 0.4398646581627901
$ws.Range("G6").Value = 0.002381189537442518
$ws.Range("I6").Value = 0.3248244616809615
$ws.Range("M6").Value = 0.7079695435168389
$ws.Range("N6").Value = 0.9329237241476882
$ws.Range("O6").Value = 1.407437291461093
$ws.Range("B7").Value = 0.3855278357625025
$ws.Range("C7").Value = 0.1152899210241856
$ws.Range("D7").Value = 0.01129876495581073
$ws.Range("F7").Value = 0.4419538898007644
$ws.Range("G7").Value = 0.002380406298180452
$ws.Range("I7").Value = 0.3240012600603563
$ws.Range("M7").Value = 0.7360179956710482
$ws.Range("N7").Value = 0.9290639380208496
$ws.Range("O7").Value = 1.410706746131268
$ws.Range("B8").Value = 0.4638371650086981
$ws.Range("C8").Value = 0.1339997442496212
$ws.Range("D8").Value = 0.01333821764798415
$ws.Range("F8").Value = 0.4519381550351653
$ws.Range("G8").Value = 0.002377131816653277
$ws.Range("I8").Value = 0.3209019774475657
$ws.Range("M8").Value = 0.8608916449777979
$ws.Range("N8").Value = 0.912980663935862
$ws.Range("O8").Value = 1.427574929379489
$ws.Range("B9").Value = 0.6164249692042176
$ws.Range("C9").Value = 0.1702043889849278
$ws.Range("D9").Value = 0.01729274212684118
$ws.Range("F9").Value = 0.4742780422288675
$ws.Range("G9").Value = 0.00237135957428936
$ws.Range("I9").Value = 0.3167510665232705
$ws.Range("M9").Value = 1.11013826240783
$ws.Range("N9").Value = 0.8848548814462553
$ws.Range("O9").Value = 1.469469044063914
$ws.Range("B10").Value = 0.7278843016820247
$ws.Range("C10").Value = 0.1964967278833001
$ws.Range("D10").Value = 0.02016956157279282
$ws.Range("F10").Value = 0.4923881525519676
$ws.Range("G10").Value = 0.002367510230535948
$ws.Range("I10").Value = 0.314885654073958
$ws.Range("M10").Value = 1.296335015977732
$ws.Range("N10").Value = 0.8662734458529044
$ws.Range("O10").Value = 1.505670680128276
$ws.Range("B11").Value = 0.7784392349846598
$ws.Range("C11").Value = 0.2083887982885813
$ws.Range("D11").Value = 0.0214718494548265
$ws.Range("F11").Value = 0.5009985520055409
$ws.Range("G11").Value = 0.002365843198813322
$ws.Range("I11").Value = 0.3142962896339832
$ws.Range("M11").Value = 1.381800810346746
$ws.Range("N11").Value = 0.8582729303225918
$ws.Range("O11").Value = 1.523327744402053
$ws.Range("B12").Value = 0.7975606165870204
$ws.Range("C12").Value = 0.212881925682467
$ws.Range("D12").Value = 0.02196404588302414
$ws.Range("F12").Value = 0.5043127836764256
$ws.Range("G12").Value = 0.00236522395746524
$ws.Range("I12").Value = 0.3141105561585285
$ws.Range("M12").Value = 1.414281721414
$ws.Range("N12").Value = 0.8553084482594571
$ws.Range("O12").Value = 1.530185701732194
$ws.Range("B13").Value = 0.793443514455987
$ws.Range("C13").Value = 0.2119147061217745
$ws.Range("D13").Value = 0.02185808542375867
$ws.Range("F13").Value = 0.5035966154167539
$ws.Range("G13").Value = 0.002365356788201202
$ws.Range("I13").Value = 0.3141488890428477
$ws.Range("M13").Value = 1.407281064679594
$ws.Range("N13").Value = 0.8559440037040034
$ws.Range("O13").Value = 1.528701074847987
$ws.Range("B14").Value = 0.7800128245672227
$ws.Range("C14").Value = 0.2087586555500138
$ws.Range("D14").Value = 0.021512362011741
$ws.Range("F14").Value = 0.5012701394675219
$ws.Range("G14").Value = 0.002365792012625151
$ws.Range("I14").Value = 0.3142802577767725
$ws.Range("M14").Value = 1.38447065139502
$ws.Range("N14").Value = 0.8580277344940086
$ws.Range("O14").Value = 1.523888509521839
$ws.Range("B15").Value = 0.7717831470480405
$ws.Range("C15").Value = 0.2068241538308371
$ws.Range("D15").Value = 0.02130047152762415
$ws.Range("F15").Value = 0.4998520987974047
$ws.Range("G15").Value = 0.002366060165041674
$ws.Range("I15").Value = 0.3143656065663478
$ws.Range("M15").Value = 1.370514045698854
$ws.Range("N15").Value = 0.8593125660229646
$ws.Range("O15").Value = 1.520963042876218
$ws.Range("B16").Value = 0.7245773169998984
$ws.Range("C16").Value = 0.1957181502399408
$ws.Range("D16").Value = 0.02008432259288639
$ws.Range("F16").Value = 0.4918329427316337
$ws.Range("G16").Value = 0.002367620860976023
$ws.Range("I16").Value = 0.3149294011673689
$ws.Range("M16").Value = 1.2907655935891
$ws.Range("N16").Value = 0.8668054079338106
$ws.Range("O16").Value = 1.504540727416014
$ws.Range("B17").Value = 0.6955790988325816
$ws.Range("C17").Value = 0.1888872333220775
$ws.Range("D17").Value = 0.0193365943545345
$ws.Range("F17").Value = 0.4870088503103744
$ws.Range("G17").Value = 0.002368599783290924
$ws.Range("I17").Value = 0.3153417846627384
$ws.Range("M17").Value = 1.242043103577814
$ws.Range("N17").Value = 0.8715179332349905
$ws.Range("O17").Value = 1.494771113534796
$ws.Range("B18").Value = 0.6788861945111648
$ws.Range("C18").Value = 0.1849518459756894
$ws.Range("D18").Value = 0.01890592118078871
$ws.Range("F18").Value = 0.4842691663727408
$ws.Range("G18").Value = 0.002369170748556487
$ws.Range("I18").Value = 0.3156033656411203
$ws.Range("M18").Value = 1.21409069316411
$ws.Range("N18").Value = 0.8742710166644763
$ws.Range("O18").Value = 1.489263741626644
$ws.Range("B19").Value = 0.6732319175969792
$ws.Range("C19").Value = 0.1836182971152027
$ws.Range("D19").Value = 0.01876000064437733
$ws.Range("F19").Value = 0.4833475638403399
$ws.Range("G19").Value = 0.002369365429076638
$ws.Range("I19").Value = 0.3156961160975769
$ws.Range("M19").Value = 1.204638552763598
$ws.Range("N19").Value = 0.8752104706749648
$ws.Range("O19").Value = 1.48741823038992
$ws.Range("B20").Value = 0.698667456455496
$ws.Range("C20").Value = 0.189615063177115
$ws.Range("D20").Value = 0.01941625360096566
$ws.Range("F20").Value = 0.487518758961258
$ws.Range("G20").Value = 0.002368494756774654
$ws.Range("I20").Value = 0.315295360345619
$ws.Range("M20").Value = 1.247222247957168
$ws.Range("N20").Value = 0.8710118709820662
$ws.Range("O20").Value = 1.49579952455997
$ws.Range("B21").Value = 0.7839583696763839
$ws.Range("C21").Value = 0.2096859417323458
$ws.Range("D21").Value = 0.02161393544827206
$ws.Range("F21").Value = 0.5019520242153845
$ws.Range("G21").Value = 0.00236566385071591
$ws.Range("I21").Value = 0.3142406539593594
$ws.Range("M21").Value = 1.391167394124111
$ws.Range("N21").Value = 0.8574139235291653
$ws.Range("O21").Value = 1.525297413830856
$ws.Range("B22").Value = 0.839568134043077
$ws.Range("C22").Value = 0.2227441691662477
$ws.Range("D22").Value = 0.02304468501012025
$ws.Range("F22").Value = 0.5116978466661521
$ws.Range("G22").Value = 0.002363883766685651
$ws.Range("I22").Value = 0.3137696860725221
$ws.Range("M22").Value = 1.485928457245251
$ws.Range("N22").Value = 0.8489066139371459
$ws.Range("O22").Value = 1.54557657672882
$ws.Range("B23").Value = 0.8099007426288267
$ws.Range("C23").Value = 0.2157802679108443
$ws.Range("D23").Value = 0.02228158680446768
$ws.Range("F23").Value = 0.5064676378249118
$ws.Range("G23").Value = 0.002364827438759948
$ws.Range("I23").Value = 0.3140010155811588
$ws.Range("M23").Value = 1.435287690455937
$ws.Range("N23").Value = 0.8534123437238819
$ws.Range("O23").Value = 1.534661430953605
$ws.Range("B24").Value = 0.6972712766619793
$ws.Range("C24").Value = 0.1892860368243987
$ws.Range("D24").Value = 0.01938024212746825
$ws.Range("F24").Value = 0.4872881241455502
$ws.Range("G24").Value = 0.002368542213662648
$ws.Range("I24").Value = 0.3153162724920335
$ws.Range("M24").Value = 1.244880574168548
$ws.Range("N24").Value = 0.8712405253388198
$ws.Range("O24").Value = 1.495334239511806
$ws.Range("B25").Value = 0.5752558107728873
$ws.Range("C25").Value = 0.1604632000748722
$ws.Range("D25").Value = 0.01622786619230965
$ws.Range("F25").Value = 0.4679375078602561
$ws.Range("G25").Value = 0.002372852066457381
$ws.Range("I25").Value = 0.3176667281389491
$ws.Range("M25").Value = 1.042200682754412
$ws.Range("N25").Value = 0.8920980690611096
$ws.Range("O25").Value = 1.457186864561322
